$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A:A").Insert()
$ws.Range("1:3").Insert()
